$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint")

# Row 11: task moved to "Done", logging 1 effort unit for day 4 (column I)
$ws.Range("F11").Value = "Done"
$ws.Range("I11").Value = 1

# Row 12: assignee added, task moved to "Done", logging 1 effort unit for day 4 (column I)
$ws.Range("D12").Value = "Predrag Dimitrijević"
$ws.Range("F12").Value = "Done"
$ws.Range("I12").Value = 1

# Row 13: assignee added, task moved to "In progress"
$ws.Range("D13").Value = "Predrag Dimitrijević"
$ws.Range("F13").Value = "In progress"

# Recalculate the workbook so dependent totals / burndown formulas refresh
$excel.Calculate()

# Update the active selection to reflect where the author was last working
$ws.Activate()
[void]$ws.Range("I13").Select()
